$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: make bold (A1:E1) ---
$ws.Range("A1:E1").Font.Bold = $true

# --- New data rows 30-34 (Line "EFOE D") ---
$ws.Range("A30").Value = "EFOE D"
$ws.Range("B30").Value = "HARTOE13"
$ws.Range("C30").Value = 5
$ws.Range("D30").Formula = "=372*497"
$ws.Range("E30").Formula = "=C30/(D30/10000)"
$ws.Range("F30").Formula = "=AVERAGE(E30:E34)"

$ws.Range("A31").Value = "EFOE D"
$ws.Range("B31").Value = "HARTOE17"
$ws.Range("C31").Value = 0
$ws.Range("D31").Formula = "=306*229"
$ws.Range("E31").Formula = "=C31/(D31/10000)"

$ws.Range("A32").Value = "EFOE D"
$ws.Range("B32").Value = "HARTOE19"
$ws.Range("C32").Value = 3
$ws.Range("D32").Formula = "=995*749"
$ws.Range("E32").Formula = "=C32/(D32/10000)"

$ws.Range("A33").Value = "EFOE D"
$ws.Range("B33").Value = "HARTOED2"
$ws.Range("C33").Value = 1
$ws.Range("D33").Formula = "=384*288"
$ws.Range("E33").Formula = "=C33/(D33/10000)"

$ws.Range("A34").Value = "EFOE D"
$ws.Range("B34").Value = "HARTOED3"
$ws.Range("C34").Value = 14
$ws.Range("D34").Formula = "=369*493"
$ws.Range("E34").Formula = "=C34/(D34/10000)"

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- View state: scroll + selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("F30").Select()

Write-Host "done"
